# case with 380 kV done
# Update simulated power-flow results (pl_mw) for rows 2-25,
# columns C,D,E,F,G,L,N,O.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.4113150465963145
$ws.Range("D2").Value = 0.04122572234083322
$ws.Range("E2").Value = 0.1738898040102441
$ws.Range("F2").Value = 1.595332469086912
$ws.Range("G2").Value = 0.002439077633301376
$ws.Range("L2").Value = 0.1518719258989734
$ws.Range("N2").Value = 2.873581687696401
$ws.Range("O2").Value = 5.764297367248616
$ws.Range("C3").Value = 0.4025426414731328
$ws.Range("D3").Value = 0.04097521872294152
$ws.Range("E3").Value = 0.1697737840431728
$ws.Range("F3").Value = 1.528102653615178
$ws.Range("G3").Value = 0.002445006840919639
$ws.Range("L3").Value = 0.1479145768503898
$ws.Range("N3").Value = 2.562605684679681
$ws.Range("O3").Value = 5.534998599706682
$ws.Range("C4").Value = 0.3974091552244943
$ws.Range("D4").Value = 0.04084533152759917
$ws.Range("E4").Value = 0.1673565443793592
$ws.Range("F4").Value = 1.487888719256986
$ws.Range("G4").Value = 0.0024488389505038
$ws.Range("L4").Value = 0.1455834063615171
$ws.Range("N4").Value = 2.371325805375818
$ws.Range("O4").Value = 5.398053472804634
$ws.Range("C5").Value = 0.395380456662167
$ws.Range("D5").Value = 0.04079838412101111
$ws.Range("E5").Value = 0.1663989955515106
$ws.Range("F5").Value = 1.471766925594494
$ws.Range("G5").Value = 0.002450448910948996
$ws.Range("L5").Value = 0.1446580720819739
$ws.Range("N5").Value = 2.293303068607429
$ws.Range("O5").Value = 5.343206354397012
$ws.Range("C6").Value = 0.3950474027081441
$ws.Range("D6").Value = 0.04079094892883006
$ws.Range("E6").Value = 0.1662416511366303
$ws.Range("F6").Value = 1.469105898554105
$ws.Range("G6").Value = 0.002450719169078277
$ws.Range("L6").Value = 0.1445059040533607
$ws.Range("N6").Value = 2.280343261403573
$ws.Range("O6").Value = 5.33415673610159
$ws.Range("C7").Value = 0.3973815398671832
$ws.Range("D7").Value = 0.04084467419834326
$ws.Range("E7").Value = 0.1673435194300446
$ws.Range("F7").Value = 1.487670222138831
$ws.Range("G7").Value = 0.002448860466967142
$ws.Range("L7").Value = 0.1455708274514933
$ws.Range("N7").Value = 2.370273851392596
$ws.Range("O7").Value = 5.397309912077901
$ws.Range("C8").Value = 0.4082376467131894
$ws.Range("D8").Value = 0.04113436324431063
$ws.Range("E8").Value = 0.1724476636962464
$ws.Range("F8").Value = 1.571929136775665
$ws.Range("G8").Value = 0.002441082374678493
$ws.Range("L8").Value = 0.1504868499159144
$ws.Range("N8").Value = 2.766433886209825
$ws.Range("O8").Value = 5.684432272078936
$ws.Range("C9").Value = 0.4315490945530769
$ws.Range("D9").Value = 0.0418938558584685
$ws.Range("E9").Value = 0.1833381646391032
$ws.Range("F9").Value = 1.745725099980206
$ws.Range("G9").Value = 0.002427341336521415
$ws.Range("L9").Value = 0.1609184648714432
$ws.Range("N9").Value = 3.540180268007646
$ws.Range("O9").Value = 6.278378879143929
$ws.Range("C10").Value = 0.4499349010714582
$ws.Range("D10").Value = 0.04257091168715732
$ws.Range("E10").Value = 0.1918894678433247
$ws.Range("F10").Value = 1.878803705533187
$ws.Range("G10").Value = 0.002418156070073038
$ws.Range("L10").Value = 0.1690778548396423
$ws.Range("N10").Value = 4.10623028343673
$ws.Range("O10").Value = 6.734186266718666
$ws.Range("C11").Value = 0.4585782732348207
$ws.Range("D11").Value = 0.04290530830628114
$ws.Range("E11").Value = 0.1959019978053576
$ws.Range("F11").Value = 1.94055330003863
$ws.Range("G11").Value = 0.002414172684371449
$ws.Range("L11").Value = 0.172900211514559
$ws.Range("N11").Value = 4.363110593465422
$ws.Range("O11").Value = 6.945899248271644
$ws.Range("C12").Value = 0.4618919340935577
$ws.Range("D12").Value = 0.04303577475471343
$ws.Range("E12").Value = 0.1974392716318931
$ws.Range("F12").Value = 1.964113322163229
$ws.Range("G12").Value = 0.002412692139293403
$ws.Range("L12").Value = 0.1743637640832247
$ws.Range("N12").Value = 4.460285735714251
$ws.Range("O12").Value = 7.026707033160847
$ws.Range("C13").Value = 0.4611764656201274
$ws.Range("D13").Value = 0.0430075051174228
$ws.Range("E13").Value = 0.1971073966083026
$ws.Range("F13").Value = 1.959031345373745
$ws.Range("G13").Value = 0.002413009764156291
$ws.Range("L13").Value = 0.1740478421134668
$ws.Range("N13").Value = 4.439361943450422
$ws.Range("O13").Value = 7.009275167049168
$ws.Range("C14").Value = 0.4588500740291863
$ws.Range("D14").Value = 0.04291596470560677
$ws.Range("E14").Value = 0.1960281120293459
$ws.Range("F14").Value = 1.942488039117279
$ws.Range("G14").Value = 0.002414050321278413
$ws.Range("L14").Value = 0.1730202948356521
$ws.Range("N14").Value = 4.371107314139522
$ws.Range("O14").Value = 6.952534540460306
$ws.Range("C15").Value = 0.4574303918063265
$ws.Range("D15").Value = 0.04286039458800417
$ws.Range("E15").Value = 0.1953693452286061
$ws.Range("F15").Value = 1.932377887238232
$ws.Range("G15").Value = 0.002414691318158281
$ws.Range("L15").Value = 0.1723929964418431
$ws.Range("N15").Value = 4.329286057409945
$ws.Range("O15").Value = 6.917862456505418
$ws.Range("C16").Value = 0.4493756962205566
$ws.Range("D16").Value = 0.04254959283818494
$ws.Range("E16").Value = 0.1916297216248282
$ws.Range("F16").Value = 1.874792821859501
$ws.Range("G16").Value = 0.00241842030380357
$ws.Range("L16").Value = 0.1688302978955676
$ws.Range("N16").Value = 4.089429168003562
$ws.Range("O16").Value = 6.720438916590751
$ws.Range("C17").Value = 0.4445062899972072
$ws.Range("D17").Value = 0.04236571691788527
$ws.Range("E17").Value = 0.1893671001002346
$ws.Range("F17").Value = 1.839778358249077
$ws.Range("G17").Value = 0.002420757748973889
$ws.Range("L17").Value = 0.1666731710256926
$ws.Range("N17").Value = 3.94211849063862
$ws.Range("O17").Value = 6.600450126405121
$ws.Range("C18").Value = 0.4417318269260591
$ws.Range("D18").Value = 0.04226243875648095
$ws.Range("E18").Value = 0.1880772180152661
$ws.Range("F18").Value = 1.819752836722699
$ws.Range("G18").Value = 0.002422120553716031
$ws.Range("L18").Value = 0.1654428473608789
$ws.Range("N18").Value = 3.857331695637754
$ws.Range("O18").Value = 6.531845942678046
$ws.Range("C19").Value = 0.4407969447700566
$ws.Range("D19").Value = 0.04222789565883289
$ws.Range("E19").Value = 0.1876424582821627
$ws.Range("F19").Value = 1.812992018896381
$ws.Range("G19").Value = 0.002422585135566624
$ws.Range("L19").Value = 0.1650280607486394
$ws.Range("N19").Value = 3.828614786364199
$ws.Range("O19").Value = 6.508687920264606
$ws.Range("C20").Value = 0.4450219227514083
$ws.Range("D20").Value = 0.04238503360718937
$ws.Range("E20").Value = 0.189606766517997
$ws.Range("F20").Value = 1.843493903818938
$ws.Range("G20").Value = 0.002420507023993686
$ws.Range("L20").Value = 0.1669017231800751
$ws.Range("N20").Value = 3.95780600327754
$ws.Range("O20").Value = 6.61318061865444
$ws.Range("C21").Value = 0.459532286431056
$ws.Range("D21").Value = 0.04294274785392105
$ws.Range("E21").Value = 0.1963446389377665
$ws.Range("F21").Value = 1.947342390297081
$ws.Range("G21").Value = 0.002413743929176309
$ws.Range("L21").Value = 0.1733216715419843
$ws.Range("N21").Value = 4.391158149571083
$ws.Range("O21").Value = 6.969183284980772
$ws.Range("C22").Value = 0.4692525077703635
$ws.Range("D22").Value = 0.04332963540378643
$ws.Range("E22").Value = 0.2008521682773079
$ws.Range("F22").Value = 2.01624508577035
$ws.Range("G22").Value = 0.002409486263784601
$ws.Range("L22").Value = 0.1776114721325968
$ws.Range("N22").Value = 4.67379181795809
$ws.Range("O22").Value = 7.205566999775442
$ws.Range("C23").Value = 0.4640428353817185
$ws.Range("D23").Value = 0.04312108369613554
$ws.Range("E23").Value = 0.1984368371623333
$ws.Range("F23").Value = 1.979375152919317
$ws.Range("G23").Value = 0.00241174385560077
$ws.Range("L23").Value = 0.1753132561374002
$ws.Range("N23").Value = 4.523002190001307
$ws.Range("O23").Value = 7.079061535084406
$ws.Range("C24").Value = 0.4447887272712023
$ws.Range("D24").Value = 0.04237629295153766
$ws.Range("E24").Value = 0.1894983793037923
$ws.Range("F24").Value = 1.841813779557327
$ws.Range("G24").Value = 0.002420620317637123
$ws.Range("L24").Value = 0.166798364125242
$ws.Range("N24").Value = 3.950713976768498
$ws.Range("O24").Value = 6.607423983390333
$ws.Range("C25").Value = 0.4250235504886177
$ws.Range("D25").Value = 0.04166765987362453
$ws.Range("E25").Value = 0.1802963383349123
$ws.Range("F25").Value = 1.697774433799424
$ws.Range("G25").Value = 0.002430897980181944
$ws.Range("L25").Value = 0.1580104031885554
$ws.Range("N25").Value = 3.331249627311138
$ws.Range("O25").Value = 5.343206354397012
